{"js": "// Replace the computed three-digit \u00f7 one-digit division answers in the\n// table cells with the new values from the commit.\n//\n// Each old value is unique within the document, so a plain text search\n// and replace (matchCase, no wildcards) on context.document.body is\n// sufficient and safe \u2014 it cannot cross-match another cell.\nconst replacements = [\n  [\"448\u00f76=74, 4\", \"384\u00f73=128, 0\"],\n  [\"952\u00f79=105, 7\", \"646\u00f77=92, 2\"],\n  [\"607\u00f76=101, 1\", \"604\u00f78=75, 4\"],\n  [\"536\u00f72=268, 0\", \"140\u00f78=17, 4\"],\n  [\"955\u00f76=159, 1\", \"836\u00f72=418, 0\"],\n  [\"643\u00f75=128, 3\", \"540\u00f79=60, 0\"],\n  [\"879\u00f76=146, 3\", \"531\u00f79=59, 0\"],\n  [\"609\u00f75=121, 4\", \"170\u00f77=24, 2\"],\n  [\"403\u00f78=50, 3\", \"168\u00f76=28, 0\"],\n  [\"144\u00f73=48, 0\", \"332\u00f73=110, 2\"],\n  [\"216\u00f77=30, 6\", \"538\u00f75=107, 3\"],\n  [\"908\u00f74=227, 0\", \"540\u00f75=108, 0\"],\n  [\"275\u00f72=137, 1\", \"784\u00f79=87, 1\"],\n  [\"833\u00f79=92, 5\", \"649\u00f78=81, 1\"],\n  [\"754\u00f76=125, 4\", \"999\u00f79=111, 0\"],\n  [\"187\u00f79=20, 7\", \"551\u00f74=137, 3\"],\n  [\"947\u00f76=157, 5\", \"721\u00f78=90, 1\"],\n  [\"162\u00f77=23, 1\", \"208\u00f72=104, 0\"],\n  [\"448\u00f74=112, 0\", \"657\u00f77=93, 6\"],\n  [\"838\u00f74=209, 2\", \"954\u00f74=238, 2\"],\n  [\"903\u00f79=100, 3\", \"579\u00f73=193, 0\"],\n  [\"163\u00f76=27, 1\", \"480\u00f78=60, 0\"],\n  [\"869\u00f74=217, 1\", \"306\u00f73=102, 0\"],\n  [\"708\u00f73=236, 0\", \"814\u00f74=203, 2\"],\n  [\"718\u00f74=179, 2\", \"192\u00f74=48, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the computed three-digit \u00f7 one-digit division answers in the\n# table cells with the new values from the commit.\n#\n# Each old value is unique within the document, so Find/Replace (exact\n# match, no wildcards) across the whole document body is safe and will\n# not cross-match any other cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"448\u00f76=74, 4\", \"384\u00f73=128, 0\"),\n    @(\"952\u00f79=105, 7\", \"646\u00f77=92, 2\"),\n    @(\"607\u00f76=101, 1\", \"604\u00f78=75, 4\"),\n    @(\"536\u00f72=268, 0\", \"140\u00f78=17, 4\"),\n    @(\"955\u00f76=159, 1\", \"836\u00f72=418, 0\"),\n    @(\"643\u00f75=128, 3\", \"540\u00f79=60, 0\"),\n    @(\"879\u00f76=146, 3\", \"531\u00f79=59, 0\"),\n    @(\"609\u00f75=121, 4\", \"170\u00f77=24, 2\"),\n    @(\"403\u00f78=50, 3\", \"168\u00f76=28, 0\"),\n    @(\"144\u00f73=48, 0\", \"332\u00f73=110, 2\"),\n    @(\"216\u00f77=30, 6\", \"538\u00f75=107, 3\"),\n    @(\"908\u00f74=227, 0\", \"540\u00f75=108, 0\"),\n    @(\"275\u00f72=137, 1\", \"784\u00f79=87, 1\"),\n    @(\"833\u00f79=92, 5\", \"649\u00f78=81, 1\"),\n    @(\"754\u00f76=125, 4\", \"999\u00f79=111, 0\"),\n    @(\"187\u00f79=20, 7\", \"551\u00f74=137, 3\"),\n    @(\"947\u00f76=157, 5\", \"721\u00f78=90, 1\"),\n    @(\"162\u00f77=23, 1\", \"208\u00f72=104, 0\"),\n    @(\"448\u00f74=112, 0\", \"657\u00f77=93, 6\"),\n    @(\"838\u00f74=209, 2\", \"954\u00f74=238, 2\"),\n    @(\"903\u00f79=100, 3\", \"579\u00f73=193, 0\"),\n    @(\"163\u00f76=27, 1\", \"480\u00f78=60, 0\"),\n    @(\"869\u00f74=217, 1\", \"306\u00f73=102, 0\"),\n    @(\"708\u00f73=236, 0\", \"814\u00f74=203, 2\"),\n    @(\"718\u00f74=179, 2\", \"192\u00f74=48, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $old,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $new,\n        2\n    )\n}\n"}
